$d = $word.ActiveDocument

function Get-ParagraphAtStart($startPos) {
  for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Start -eq $startPos) {
      return $pp
    }
  }
  return $null
}

# -----------------------------------------------------------------
# Change 1: split the run "可以使用" into three runs "可以" / "基本" / "使用"
# (all three keep the same rPr, including rFonts hint="eastAsia")
# -----------------------------------------------------------------
$find1 = $d.Content.Find
$find1.Execute("可以使用")
$r1 = $find1.Parent
$s1 = $r1.Start

# Remove the whole original run/text...
$r1.Delete()

# ...then insert the three-run replacement right at the same spot.
$ins1 = $d.Range($s1, $s1)
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="微软雅黑" w:hAnsi="微软雅黑" w:hint="eastAsia"/><w:color w:val="414141"/><w:szCs w:val="21"/></w:rPr><w:t>可以</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微软雅黑" w:hAnsi="微软雅黑" w:hint="eastAsia"/><w:color w:val="414141"/><w:szCs w:val="21"/></w:rPr><w:t>基本</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微软雅黑" w:hAnsi="微软雅黑" w:hint="eastAsia"/><w:color w:val="414141"/><w:szCs w:val="21"/></w:rPr><w:t>使用</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$ins1.InsertXML($xml1)

# -----------------------------------------------------------------
# Change 2: paragraph "参与FlinkCDC的部署和使用。" — drop hint="eastAsia"
# from the paragraph mark's run properties (pPr/rPr/rFonts).
# -----------------------------------------------------------------
$find2 = $d.Content.Find
$find2.Execute("参与FlinkCDC的部署和使用。")
$r2 = $find2.Parent
$p2 = Get-ParagraphAtStart $r2.Start
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="092430B4" w14:textId="3DE854B8" w:rsidR="00C6011E" w:rsidRPr="00C6011E" w:rsidRDefault="00C6011E" w:rsidP="00C6011E"><w:pPr><w:pStyle w:val="a7"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="20"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="3828"/><w:tab w:val="left" w:pos="7655"/></w:tabs><w:adjustRightInd w:val="0"/><w:snapToGrid w:val="0"/><w:ind w:firstLineChars="0"/><w:rPr><w:rFonts w:ascii="微软雅黑" w:hAnsi="微软雅黑"/><w:color w:val="414141"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="微软雅黑" w:hAnsi="微软雅黑" w:hint="eastAsia"/><w:color w:val="414141"/><w:szCs w:val="21"/></w:rPr><w:t>参与FlinkCDC的部署和使用。</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p2.Range.InsertXML($xml2)

# -----------------------------------------------------------------
# Change 3: paragraph "共学营期间设计基于ERC721的NFT程序。" — drop
# hint="eastAsia" from the paragraph mark's run properties (pPr/rPr/rFonts).
# -----------------------------------------------------------------
$find3 = $d.Content.Find
$find3.Execute("共学营期间设计基于ERC721")
$r3 = $find3.Parent
$p3 = Get-ParagraphAtStart $r3.Start
$xml3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="62A4A501" w14:textId="72D454D5" w:rsidR="00A33EAB" w:rsidRPr="00070FAE" w:rsidRDefault="00A33EAB" w:rsidP="009E16C4"><w:pPr><w:pStyle w:val="a7"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="20"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="3828"/><w:tab w:val="left" w:pos="7655"/></w:tabs><w:adjustRightInd w:val="0"/><w:snapToGrid w:val="0"/><w:ind w:firstLineChars="0"/><w:rPr><w:rFonts w:ascii="微软雅黑" w:hAnsi="微软雅黑"/><w:color w:val="414141"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="微软雅黑" w:hAnsi="微软雅黑" w:hint="eastAsia"/><w:color w:val="414141"/><w:szCs w:val="21"/></w:rPr><w:t>共学营期间设计基于ERC721</w:t></w:r><w:r w:rsidR="0095639C"><w:rPr><w:rFonts w:ascii="微软雅黑" w:hAnsi="微软雅黑" w:hint="eastAsia"/><w:color w:val="414141"/><w:szCs w:val="21"/></w:rPr><w:t>的</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微软雅黑" w:hAnsi="微软雅黑" w:hint="eastAsia"/><w:color w:val="414141"/><w:szCs w:val="21"/></w:rPr><w:t>NFT</w:t></w:r><w:r w:rsidR="00091F1A"><w:rPr><w:rFonts w:ascii="微软雅黑" w:hAnsi="微软雅黑" w:hint="eastAsia"/><w:color w:val="414141"/><w:szCs w:val="21"/></w:rPr><w:t>程序</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="微软雅黑" w:hAnsi="微软雅黑" w:hint="eastAsia"/><w:color w:val="414141"/><w:szCs w:val="21"/></w:rPr><w:t>。</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p3.Range.InsertXML($xml3)

Write-Host "done"
